# Edit sur_tracking.tpl.xlsx per the "modify request from bureau" commit.
#
# Summary of the change:
#  - Drop the old standalone title row ("桃園市ＯＯ地政事務所"); the remaining
#    title ("測量案件管制清冊") becomes the new (single) title row.
#  - The old header row becomes the new row 2 and gains a 10th column (J)
#    labelled "核章欄" ("sign-off column").
#  - Column headers are renamed/reordered (收件年/收件字/收件號/複丈原因/
#    收件日期/複丈日期/逾期日期/測量員/處理情形/核章欄).
#  - Column widths are resized and a new column J is added.
#  - The selection/active cell moves to J16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: remove the old title row (row 1). Because the remaining title
# row / header row / data rows were merged (A1:I1) and (A2:I2) before,
# Excel automatically collapses the merge/shift so that:
#   old row2 (title "測量案件管制清冊")   -> new row1
#   old row3 (headers)                   -> new row2
#   old rows4-71 (data)                  -> new rows3-70
# ---------------------------------------------------------------------
$ws.Rows.Item(1).Delete()

# ---------------------------------------------------------------------
# Step 2: new row 1 - title row.
# A1 keeps the 標楷體/size18 "title" look with a bottom border and
# centered text; B1:J1 use the plain default font with the same bottom
# border (J1 keeps default/vertical-only alignment, matching a normal
# unlabeled cell).
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "測量案件管制清冊"
$ws.Range("A1").Borders.Item(9).LineStyle = 1
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108

# Reset B1:J1 back to the workbook's plain default formatting (fresh,
# never-styled cell) by pasting formats from a pristine cell on the
# blank helper sheet, then layer the new border/alignment on top.
$pristine = $wb.Worksheets.Item(2)
$pristine.Cells.Item(1, 1).Copy()
$ws.Range("B1:J1").PasteSpecial(-4122)

$ws.Range("B1:I1").Borders.Item(9).LineStyle = 1
$ws.Range("B1:I1").HorizontalAlignment = -4108
$ws.Range("B1:I1").VerticalAlignment = -4108
$ws.Range("J1").Borders.Item(9).LineStyle = 1
$ws.Range("J1").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Step 3: new row 2 - header row. A2:I2 already carry the former header
# styling (boxed border, centered, 標楷體). Give J2 the same look by
# copying the format from A2, then set the text for every header cell.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("J2").PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = "收件年"
$ws.Cells.Item(2, 2).Value = "收件字"
$ws.Cells.Item(2, 3).Value = "收件號"
$ws.Cells.Item(2, 4).Value = "複丈原因"
$ws.Cells.Item(2, 5).Value = "收件日期"
$ws.Cells.Item(2, 6).Value = "複丈日期"
$ws.Cells.Item(2, 7).Value = "逾期日期"
$ws.Cells.Item(2, 8).Value = "測量員"
$ws.Cells.Item(2, 9).Value = "處理情形"
$ws.Cells.Item(2, 10).Value = "核章欄"

# ---------------------------------------------------------------------
# Step 4: resize columns A-I and size the new column J.
# ---------------------------------------------------------------------
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 9.857142857142858
$ws.Range("D1").EntireColumn.ColumnWidth = 14.857142857142858
$ws.Range("E1:G1").EntireColumn.ColumnWidth = 11.857142857142858
$ws.Range("H1").EntireColumn.ColumnWidth = 19.857142857142858
$ws.Range("I1").EntireColumn.ColumnWidth = 29.857142857142858
$ws.Range("J1").EntireColumn.ColumnWidth = 18.857142857142858

# ---------------------------------------------------------------------
# Step 5: move the active selection to J16, matching the saved view.
# ---------------------------------------------------------------------
$ws.Range("J16").Select() | Out-Null
